# Scheduled-runner market data refresh for the Leve profit tracker workbook.
# Recomputed currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ /
# LevePriceHQ / LeveProfitNQ / LeveProfitHQ (columns H:N) for the leve rows whose market-board
# data changed, across all eight job sheets. A handful of rows flip which HQ/NQ profit column
# applies (the craft goes from NQ-only to HQ-viable or vice versa), so those cells are added or
# cleared outright instead of merely overwritten, matching upstream.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 241
$ws.Range("I2").Value = 241
$ws.Range("K2").Value = 241
$ws.Range("M2").Value = -128
$ws.Range("H129").Value = 825
$ws.Range("I129").Value = 487.57693
$ws.Range("J129").Value = 1622.5454
$ws.Range("K129").Value = 1462.73079
$ws.Range("L129").Value = 4867.6362
$ws.Range("M129").Value = 3537.26921
$ws.Range("N129").Value = -14867.6362
$ws.Range("H132").Value = 1775.8933
$ws.Range("I132").Value = 1602.459
$ws.Range("J132").Value = 2531.5715
$ws.Range("K132").Value = 4807.377
$ws.Range("L132").Value = 7594.7145
$ws.Range("M132").Value = -2277.377
$ws.Range("N132").Value = -12654.7145
$ws.Range("H137").Value = 3572769
$ws.Range("I137").Value = 1614187
$ws.Range("J137").Value = 9092409
$ws.Range("K137").Value = 4842561
$ws.Range("L137").Value = 27277227
$ws.Range("M137").Value = -4840011
$ws.Range("N137").Value = -27282327
$ws.Range("H141").Value = 2605717
$ws.Range("I141").Value = 1288.3265
$ws.Range("J141").Value = 11113517
$ws.Range("K141").Value = 3864.979499999999
$ws.Range("L141").Value = 33340551
$ws.Range("M141").Value = 1315.020500000001
$ws.Range("N141").Value = -33350911

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 171248.5
$ws.Range("J24").Value = 171248.5
$ws.Range("L24").Value = 171248.5
$ws.Range("N24").Value = -171996.5
$ws.Range("H32").Value = 2649581
$ws.Range("I32").Value = 5001.357
$ws.Range("J32").Value = 22844554
$ws.Range("K32").Value = 5001.357
$ws.Range("L32").Value = 22844554
$ws.Range("M32").Value = -4714.357
$ws.Range("N32").Value = -22845128
$ws.Range("H61").Value = 1723.1786
$ws.Range("I61").Value = 1701.8846
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1701.8846
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1489.8846
$ws.Range("N61").Value = -2424
$ws.Range("H88").Value = 23594.223
$ws.Range("I88").Value = 1585
$ws.Range("J88").Value = 41201.6
$ws.Range("K88").Value = 1585
$ws.Range("L88").Value = 41201.6
$ws.Range("M88").Value = -1179
$ws.Range("N88").Value = -42013.6
$ws.Range("H91").Value = 23594.223
$ws.Range("I91").Value = 1585
$ws.Range("J91").Value = 41201.6
$ws.Range("K91").Value = 1585
$ws.Range("L91").Value = 41201.6
$ws.Range("M91").Value = -181
$ws.Range("N91").Value = -44009.6
$ws.Range("H100").Value = 171248.5
$ws.Range("J100").Value = 171248.5
$ws.Range("L100").Value = 171248.5
$ws.Range("N100").Value = -173412.5
$ws.Range("H132").Value = 91288.55499999999
$ws.Range("I132").Value = 120151.45
$ws.Range("J132").Value = 4699.857
$ws.Range("K132").Value = 360454.35
$ws.Range("L132").Value = 14099.571
$ws.Range("M132").Value = -357924.35
$ws.Range("N132").Value = -19159.571
$ws.Range("H136").Value = 1723.1786
$ws.Range("I136").Value = 1701.8846
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5105.6538
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2555.6538
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 42354.75
$ws.Range("J74").Value = 42354.75
$ws.Range("L74").Value = 42354.75
$ws.Range("N74").Value = -44226.75
$ws.Range("H77").Value = 42354.75
$ws.Range("J77").Value = 42354.75
$ws.Range("L77").Value = 127064.25
$ws.Range("N77").Value = -136424.25
$ws.Range("H86").Value = 2043.5518
$ws.Range("I86").Value = 2094.087
$ws.Range("J86").Value = 1849.8334
$ws.Range("K86").Value = 2094.087
$ws.Range("L86").Value = 1849.8334
$ws.Range("M86").Value = -971.087
$ws.Range("N86").Value = -4095.8334
$ws.Range("H89").Value = 2043.5518
$ws.Range("I89").Value = 2094.087
$ws.Range("J89").Value = 1849.8334
$ws.Range("K89").Value = 10470.435
$ws.Range("L89").Value = 9249.166999999999
$ws.Range("M89").Value = -4854.434999999999
$ws.Range("N89").Value = -20481.167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 961.9583
$ws.Range("I58").Value = 959.42224
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 959.42224
$ws.Range("L58").Value = 1000
$ws.Range("M58").Value = -756.42224
$ws.Range("N58").Value = -1406
$ws.Range("H122").Value = 1565.0526
$ws.Range("I122").Value = 687.2857
$ws.Range("J122").Value = 4022.8
$ws.Range("K122").Value = 2061.8571
$ws.Range("L122").Value = 12068.4
$ws.Range("M122").Value = 388.1428999999998
$ws.Range("N122").Value = -16968.4
$ws.Range("H134").Value = 3630.12
$ws.Range("I134").Value = 4236.364
$ws.Range("J134").Value = 2453.2942
$ws.Range("K134").Value = 12709.092
$ws.Range("L134").Value = 7359.882599999999
$ws.Range("M134").Value = -10174.092
$ws.Range("N134").Value = -12429.8826
$ws.Range("H136").Value = 961.9583
$ws.Range("I136").Value = 959.42224
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 2878.26672
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -328.2667200000001
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 2351.25
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 2351.25
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 7053.75
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -10797.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 25117.393
$ws.Range("I70").Value = 30117.072
$ws.Range("J70").Value = 4618.7
$ws.Range("K70").Value = 30117.072
$ws.Range("L70").Value = 4618.7
$ws.Range("M70").Value = -29847.072
$ws.Range("N70").Value = -5158.7
$ws.Range("H73").Value = 25117.393
$ws.Range("I73").Value = 30117.072
$ws.Range("J73").Value = 4618.7
$ws.Range("K73").Value = 30117.072
$ws.Range("L73").Value = 4618.7
$ws.Range("M73").Value = -29181.072
$ws.Range("N73").Value = -6490.7
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H132").Value = 1796.4889
$ws.Range("I132").Value = 1382.0938
$ws.Range("J132").Value = 2816.5386
$ws.Range("K132").Value = 4146.2814
$ws.Range("L132").Value = 8449.6158
$ws.Range("M132").Value = -1616.2814
$ws.Range("N132").Value = -13509.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 28369
$ws.Range("J63").Value = 28369
$ws.Range("L63").Value = 28369
$ws.Range("N63").Value = -29867
$ws.Range("H66").Value = 28369
$ws.Range("J66").Value = 28369
$ws.Range("L66").Value = 85107
$ws.Range("N66").Value = -92595
$ws.Range("H118").Value = 32021.5
$ws.Range("J118").Value = 32021.5
$ws.Range("L118").Value = 32021.5
$ws.Range("N118").Value = -35335.5
$ws.Range("H133").Value = 41162
$ws.Range("J133").Value = 41162
$ws.Range("L133").Value = 41162
$ws.Range("N133").Value = -46222
$ws.Range("H134").Value = 34376.945
$ws.Range("J134").Value = 34376.945
$ws.Range("L134").Value = 34376.945
$ws.Range("N134").Value = -44516.945
$ws.Range("H136").Value = 2504.7273
$ws.Range("I136").Value = 2405.7778
$ws.Range("J136").Value = 2950
$ws.Range("K136").Value = 7217.3334
$ws.Range("L136").Value = 8850
$ws.Range("M136").Value = -4667.3334
$ws.Range("N136").Value = -13950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 27637.857
$ws.Range("J75").Value = 27637.857
$ws.Range("L75").Value = 27637.857
$ws.Range("N75").Value = -29509.857
$ws.Range("H78").Value = 27637.857
$ws.Range("J78").Value = 27637.857
$ws.Range("L78").Value = 82913.571
$ws.Range("N78").Value = -92273.571
$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 20000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 20000
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("M88").Value = -19594
$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 20000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 20000
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("M91").Value = -18596
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H121").Value = 31420
$ws.Range("J121").Value = 31420
$ws.Range("L121").Value = 31420
$ws.Range("N121").Value = -34914
$ws.Range("H126").Value = 2092.4614
$ws.Range("I126").Value = 1382
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 4146
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -1676
$ws.Range("N126").Value = -22940
$ws.Range("H132").Value = 1606.6316
$ws.Range("I132").Value = 1889.06
$ws.Range("J132").Value = 1063.5
$ws.Range("K132").Value = 5667.18
$ws.Range("L132").Value = 3190.5
$ws.Range("M132").Value = -3137.18
$ws.Range("N132").Value = -8250.5
$ws.Range("H136").Value = 2133.3076
$ws.Range("I136").Value = 1965.75
$ws.Range("J136").Value = 2899.2856
$ws.Range("K136").Value = 5897.25
$ws.Range("L136").Value = 8697.856800000001
$ws.Range("M136").Value = -3347.25
$ws.Range("N136").Value = -13797.8568
